# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: title + link update
$ws.Range("D9").Value = "1-2과목만 듣고 싶다는 분들에게"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/1-2-courses-only/#utm_source=rss&utm_medium=rss&utm_campaign=1-2-courses-only"

# Row 26: title update
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 32: title + link update
$ws.Range("D32").Value = "pandas 컬럼값 조건 변경"
$ws.Range("E32").Value = "https://dodonam.tistory.com/361"
